$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full dataset for the "Kosten" timetable after dropping the "id" column and
# extending the series with more rows (time, Beste Werte, Optimale Werte).
$data = @(
    @(5,  33850, 15150),
    @(10, 30250, 15150),
    @(15, 28950, 15150),
    @(20, 26950, 15150),
    @(25, 25950, 15150),
    @(30, 24450, 15150),
    @(35, 16350, 15150),
    @(40, 15850, 15150),
    @(45, 15850, 15150),
    @(50, 15850, 15150),
    @(55, 15850, 15150),
    @(60, 15850, 15150),
    @(65, 15750, 15150)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Extend the chart series ranges so they cover the full A2:C14 data block.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser1 = $chart.SeriesCollection().Item(1)
$ser2 = $chart.SeriesCollection().Item(2)

$ser1.Formula = "=SERIES('Kosten'!B1,'Kosten'!`$A`$2:`$A`$14,'Kosten'!`$B`$2:`$B`$14,1)"
$ser2.Formula = "=SERIES('Kosten'!C1,'Kosten'!`$A`$2:`$A`$14,'Kosten'!`$C`$2:`$C`$14,2)"
